$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '63.798.30'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  -0.01%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.620.15'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  -0.11%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '594.95'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -0.12%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '150.87'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +0.90%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('E8').Value = '  -0.40%  '
$ws.Range('E9').Value = '  +4.48%  '
$ws.Range('B10').Value = 'Cardano'
$ws.Range('C10').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.395'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +3.54%  '
$ws.Range('B11').Value = 'Toncoin'
$ws.Range('C11').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.80'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +3.62%  '
$ws.Range('E12').Value = '  +0.97%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '27.96'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +1.51%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '3.091.41'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -0.06%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '63.648.62'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -0.02%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000163'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +9.91%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.647.57'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +1.00%  '
$ws.Range('E18').Value = '  +0.58%  '
$ws.Range('E19').Value = '  +4.03%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '348.78'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +0.04%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '7.01'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +2.15%  '
$ws.Range('E22').Value = '  +0.20%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '67.42'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +1.97%  '
$ws.Range('E24').Value = '  -2.63%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.23'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +0.35%  '
$ws.Range('E26').Value = '  -0.23%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.40'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +3.85%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '552.06'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +1.44%  '
$ws.Range('E29').Value = '  -0.87%  '
$ws.Range('E30').Value = '  -0.09%  '
$ws.Range('E31').Value = '  +1.67%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.0₃0888'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +5.20%  '
$ws.Range('E33').Value = '  +2.29%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.41'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +3.82%  '
$ws.Range('E35').Value = '  +1.94%  '
$ws.Range('B36').Value = 'PolygonEcosystemToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.418'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +2.92%  '
$ws.Range('B37').Value = 'Monero'
$ws.Range('C37').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '164.50'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -2.08%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.00'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +2.60%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '19.82'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +2.27%  '
$ws.Range('E40').Value = '  -0.01%  '
$ws.Range('E41').Value = '  -0.10%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '167.64'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -1.23%  '
$ws.Range('E43').Value = '  +4.30%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '23.70'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +10.72%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0587'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -1.04%  '
$ws.Range('E46').Value = '  +8.97%  '
$ws.Range('E47').Value = '  +1.71%  '
$ws.Range('E48').Value = '  +3.23%  '
$ws.Range('E49').Value = '  +0.15%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '19.26'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +0.36%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0₆0232'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +19.26%  '
